$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 65, pushing the existing row 65 (and everything
# below it) down by one. The former row 65 becomes row 66 unchanged.
$ws.Rows.Item(65).Insert()

# Populate the newly inserted row 65 with the new weekly price record.
$ws.Range("A65").Value = 10
$ws.Range("B65").Value = "Vega Modelo de Temuco"
$ws.Range("C65").Value = "La Araucanía"
$ws.Range("D65").Value = 44595
$ws.Range("E65").Value = 9
$ws.Range("F65").Value = 100112030
$ws.Range("G65").Value = "Poroto granado"
$ws.Range("H65").Value = "Sin especificar"
$ws.Range("I65").Value = "Primera"
$ws.Range("J65").Value = 100
$ws.Range("K65").Value = 28000
$ws.Range("L65").Value = 28000
$ws.Range("M65").Value = 28000
$ws.Range("N65").Value = "$/saco 25 kilos"
$ws.Range("O65").Value = "Región del Maule"
$ws.Range("P65").Value = 1120
$ws.Range("Q65").Value = 25
$ws.Range("R65").Value = "Hortaliza"

# Make sure the date cell keeps the date/time number format used by the
# rest of column D (style index 2 in the original workbook).
$ws.Range("D65").NumberFormat = $ws.Range("D66").NumberFormat
